# payment_request_form.xlsx tidy-up
#
# Summary of the edit (per commit message "tidy, tidy PRF template"):
#  - Drop the "Special Instructions" row entirely.
#  - Drop the "Comments" column from the Payment Details block.
#  - Re-order the top-level fields: PO Number moves up (right after
#    Program), Qualified Receiver Name / Date Payment Authorized /
#    Expense Authority Name / Account Coding move down, after the
#    Payment Details block, with Account Coding now last.
#  - The Payment Details block (header + 2 detail rows + Total Payment)
#    moves up to directly follow Supplier/Invoice info.
#  - Sheet shrinks from 14 used rows to 13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- start from a clean slate, and normalize formatting on the rows we
#     are about to repopulate so no stale per-cell style (e.g. the old
#     bold "Payment Details"/"Comments" header row) bleeds through -------
$ws.Range("A1:E13").UnMerge()
$ws.Range("A1:E13").ClearContents()
$ws.Range("A1:E13").Font.Bold = $false

# the sheet now only needs 13 rows; drop the old trailing 14th row
$ws.Rows(14).Delete()

# --- row 1: Program -------------------------------------------------------
$ws.Range("A1").Value = "Program"
$ws.Range("A1").Font.Bold = $true
$ws.Range("B1").Value = "Dormant Sites Reclamation Program"
$ws.Range("B1:E1").Merge()

# --- row 2: PO Number (moved up from old row 6) ---------------------------
$ws.Range("A2").Value = "PO Number"
$ws.Range("A2").Font.Bold = $true
$ws.Range("B2").Value = "{d.po_number}"
$ws.Range("B2:E2").Merge()

# --- row 3: Supplier Name --------------------------------------------------
$ws.Range("A3").Value = "Supplier Name"
$ws.Range("A3").Font.Bold = $true
$ws.Range("B3").Value = "{d.supplier_name}"
$ws.Range("B3:E3").Merge()

# --- row 4: Supplier Address -----------------------------------------------
$ws.Range("A4").Value = "Supplier Address"
$ws.Range("A4").Font.Bold = $true
$ws.Range("B4").Value = "{d.supplier_address}"
$ws.Range("B4:E4").Merge()

# --- row 5: Invoice Number --------------------------------------------------
$ws.Range("A5").Value = "Invoice Number"
$ws.Range("A5").Font.Bold = $true
$ws.Range("B5").Value = "{d.invoice_number}"
$ws.Range("B5:E5").Merge()

# --- row 6: Payment Details header (moved up from old row 11, no Comments) -
$ws.Range("A6").Value = "Payment Details"
$ws.Range("A6").Font.Bold = $true
$ws.Range("B6").Value = "Agreement Number"
$ws.Range("B6").Font.Bold = $true
$ws.Range("C6").Value = "Unique ID"
$ws.Range("C6").Font.Bold = $true
$ws.Range("D6").Value = "Amount"
$ws.Range("D6").Font.Bold = $true
$ws.Range("E6").Font.Bold = $true

# --- row 7: payment_details[i] (no Comments) --------------------------------
$ws.Range("A7").Value = " "
$ws.Range("A7").Font.Bold = $true
$ws.Range("B7").Value = "{d.payment_details[i].agreement_number}"
$ws.Range("C7").Value = "{d.payment_details[i].unique_id}"
$ws.Range("D7").Value = "{d.payment_details[i].amount}"

# --- row 8: payment_details[i+1] (no Comments) ------------------------------
$ws.Range("A8").Value = " "
$ws.Range("A8").Font.Bold = $true
$ws.Range("B8").Value = "{d.payment_details[i+1].agreement_number}"
$ws.Range("C8").Value = "{d.payment_details[i+1].unique_id}"
$ws.Range("D8").Value = "{d.payment_details[i+1].amount}"

# --- row 9: Total Payment ----------------------------------------------------
$ws.Range("A9").Value = "Total Payment"
$ws.Range("A9").Font.Bold = $true
$ws.Range("B9").Value = " "
$ws.Range("C9").Value = " "
$ws.Range("D9").Value = "{d.total_payment}"

# --- row 10: Qualified Receiver Name (moved down from old row 7) -----------
$ws.Range("A10").Value = "Qualified Receiver Name"
$ws.Range("A10").Font.Bold = $true
$ws.Range("B10").Value = "{d.qualified_receiver_name}"
$ws.Range("B10:E10").Merge()

# --- row 11: Date Payment Authorized (moved down from old row 8) -----------
$ws.Range("A11").Value = "Date Payment Authorized"
$ws.Range("A11").Font.Bold = $true
$ws.Range("B11").Value = "{d.date_payment_authorized}"
$ws.Range("B11:E11").Merge()
$ws.Rows(11).RowHeight = 15

# --- row 12: Expense Authority Name (moved down from old row 9) ------------
$ws.Range("A12").Value = "Expense Authority Name"
$ws.Range("A12").Font.Bold = $true
$ws.Range("B12").Value = "{d.expense_authority_name}"
$ws.Range("B12:E12").Merge()

# --- row 13: Account Coding (moved down from old row 2, now last) ----------
$ws.Range("A13").Value = "Account Coding"
$ws.Range("A13").Font.Bold = $true
$ws.Range("B13").Value = "{d.account_coding}"
$ws.Range("B13:E13").Merge()

# --- selection bookkeeping ---------------------------------------------------
$ws.Range("B20").Select()
